# feat(CWL): allow custom zone asset loading
# Adds two new localization rows (id / ja / cn) to the "General" sheet:
#   row 131: cwl_ui_export_zone  -> "CWL/ゾーンを書き出す" / "CWL/导出地图"
#   row 132: cwl_relocate_zone   -> "relocated zone > {0}:{1}\n> {2}" / "重定向地图 > {0}:{1}\n> {2}"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 131 / 132 ids + ja text (matches authoring order so the
# shared-string table is interned in the same sequence as upstream) --
$ws.Range("A131").Value = "cwl_ui_export_zone"
$ws.Range("C131").Value = "CWL/ゾーンを書き出す"
$ws.Range("A132").Value = "cwl_relocate_zone"
$ws.Range("C132").Value = "relocated zone > {0}:{1}`n> {2}"

# --- cn text, row 131 then row 132 -----------------------------------
$ws.Range("D131").Value = "CWL/导出地图"
$ws.Range("D132").Value = "重定向地图 > {0}:{1}`n> {2}"

# Row 131 already uses the same visual style as its neighbours (24),
# but row 132's C/D cells need the wrap-text variant (25) just like the
# analogous "relocated book/drama/..." rows above it. Pull that exact
# formatting over via a format-only copy/paste so the workbook keeps
# reusing the existing style records instead of minting new ones.
$ws.Range("C129").Copy() | Out-Null
$ws.Range("C132").PasteSpecial(-4122) | Out-Null
$ws.Range("D129").Copy() | Out-Null
$ws.Range("D132").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 132 is a two-line message, so it needs to grow to match the other
# wrapped rows (e.g. row 129) instead of the default single-line height.
$ws.Rows.Item(132).RowHeight = 46.5

# Leave the workbook scrolled/selected on the newly added content, same
# as the authored selection state.
$ws.Range("D132").Select() | Out-Null
